# Edit script: "Added Counts for 10 projects"
# Adds a new header row (row 33) describing the 10 new project columns
# and 7 new data rows (34-40) with the corresponding counts, reproducing
# the OOXML diff (new shared strings, 3 new cell styles, new sheet rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: header row for the 10 new project count columns ---
$ws.Range("A33").Value = "Types"
$ws.Range("A33").VerticalAlignment = -4108
$ws.Range("A33").Font.Bold = $true

$ws.Range("B33").Value = "Quartz Scheduler"
$ws.Range("B33").VerticalAlignment = -4108

$ws.Range("C33").Value = "Open EMRConect"
$ws.Range("C33").VerticalAlignment = -4108
$ws.Range("C33").WrapText = $true

$ws.Range("D33").Value = "Wallet"
$ws.Range("D33").VerticalAlignment = -4108

$ws.Range("E33").Value = "Secure Banking System"
$ws.Range("E33").VerticalAlignment = -4108
$ws.Range("E33").WrapText = $true

$ws.Range("F33").Value = "Calendar System"
$ws.Range("F33").VerticalAlignment = -4108
$ws.Range("F33").WrapText = $true

$ws.Range("G33").Value = "Time4J"
$ws.Range("G33").VerticalAlignment = -4108
$ws.Range("G33").WrapText = $true

$ws.Range("H33").Value = "Voj "
$ws.Range("H33").VerticalAlignment = -4108
$ws.Range("H33").WrapText = $true

$ws.Range("I33").Value = "Core Flight Systm(CFS) and data Dictionary(CCDD) Utility"
$ws.Range("I33").VerticalAlignment = -4108
$ws.Range("I33").WrapText = $true

$ws.Range("J33").Value = "Dert"
$ws.Range("J33").VerticalAlignment = -4108
$ws.Range("J33").WrapText = $true

$ws.Range("K33").Value = "Hyper realistic zombie"
$ws.Range("K33").VerticalAlignment = -4108
$ws.Range("K33").WrapText = $true

$ws.Rows.Item(33).RowHeight = 96

# --- Rows 34-40: the count data for each project ---
# Row 34
$ws.Range("A34").Value = "Nested type count"
$ws.Range("A34").VerticalAlignment = -4108
$ws.Range("B34").Value = 71
$ws.Range("B34").VerticalAlignment = -4108
$ws.Range("C34").Value = 22
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 20
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 1
$ws.Range("I34").Value = 4
$ws.Range("J34").Value = 14
$ws.Range("K34").Value = 0

# Row 35
$ws.Range("A35").Value = "Local type count"
$ws.Range("A35").VerticalAlignment = -4108
$ws.Range("B35").Value = 0
$ws.Range("B35").VerticalAlignment = -4108
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0

# Row 36
$ws.Range("A36").Value = "Anonymous type count"
$ws.Range("A36").VerticalAlignment = -4108
$ws.Range("B36").Value = 0
$ws.Range("B36").VerticalAlignment = -4108
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0

# Row 37
$ws.Range("A37").Value = "Other Interface/class Decl (non nest/non local/non anon)"
$ws.Range("A37").VerticalAlignment = -4108
$ws.Range("B37").Value = 459
$ws.Range("B37").VerticalAlignment = -4108
$ws.Range("C37").Value = 177
$ws.Range("D37").Value = 11
$ws.Range("E37").Value = 61
$ws.Range("F37").Value = 45
$ws.Range("G37").Value = 58
$ws.Range("H37").Value = 39
$ws.Range("I37").Value = 5
$ws.Range("J37").Value = 334
$ws.Range("K37").Value = 20

# Row 38
$ws.Range("A38").Value = "Primitive type Count"
$ws.Range("A38").VerticalAlignment = -4108
$ws.Range("B38").Value = 1181
$ws.Range("B38").VerticalAlignment = -4108
$ws.Range("C38").Value = 542
$ws.Range("D38").Value = 3
$ws.Range("E38").Value = 112
$ws.Range("F38").Value = 116
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 165
$ws.Range("I38").Value = 117
$ws.Range("J38").Value = 1088
$ws.Range("K38").Value = 79

# Row 39
$ws.Range("A39").Value = "Annotation type Count"
$ws.Range("A39").VerticalAlignment = -4108
$ws.Range("B39").Value = 517
$ws.Range("B39").VerticalAlignment = -4108
$ws.Range("C39").Value = 78
$ws.Range("D39").Value = 27
$ws.Range("E39").Value = 311
$ws.Range("F39").Value = 62
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 74
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 190
$ws.Range("K39").Value = 7

# Row 40
$ws.Range("A40").Value = "Import Declarations Count"
$ws.Range("A40").VerticalAlignment = -4108
$ws.Range("B40").Value = 0
$ws.Range("B40").VerticalAlignment = -4108
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0

# Restore the active selection as seen in the authored workbook
$ws.Range("B46").Select()
